# Rewrites the article's opening title block into pandoc-style
# "Title" / "Authors" paragraphs:
#   - "Thank You!"    (Heading1, wrapped in a "thank-you" bookmark)
#       -> Title-styled paragraph, split into 4 runs: "Thank" / " " / "You" / "!"
#   - "By Dorothy Day" (bold run)
#       -> Authors-styled paragraph, split into 3 runs: "Dorothy" / " " / "Day"
# and drops the now-unused "thank-you" bookmark entirely.

$d = $word.ActiveDocument

# --- Step 1: clear out the old "Thank You!" (Heading1) paragraph and strip the
# "thank-you" bookmark that wraps it. Deleting the paragraph's content leaves the
# (now empty-width) bookmarkStart/bookmarkEnd pair sitting back-to-back at the very
# start of the document; two collapsed-range deletes at that point remove them.
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Delete()
$d.Range(0, 0).Delete()
$d.Range(0, 0).Delete()

# --- Step 2: insert a fresh paragraph for the title, styled "Title", ahead of the
# (now first) "By Dorothy Day" paragraph, and fill it with four discrete runs via
# raw OOXML so the run boundaries match the target exactly (plain text insertion
# would simply merge same-formatted runs together).
$d.Paragraphs.Item(1).Range.InsertParagraphBefore()
$newTitle = $d.Paragraphs.Item(1)
$titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:pStyle w:val="Title"/></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Thank</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">You</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">!</w:t></w:r>' +
  '</w:p>'
$newTitle.Range.InsertXML($titleXml) | Out-Null

# --- Step 3: replace the bold "By Dorothy Day" paragraph with an Authors-styled
# paragraph made of three runs: "Dorothy" / " " / "Day" (the "By " prefix is gone
# and so is the explicit bold run formatting -- "Authors" carries its own look).
$authorPara = $d.Paragraphs.Item(2)
$authorXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:pStyle w:val="Authors"/></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">Day</w:t></w:r>' +
  '</w:p>'
$authorPara.Range.InsertXML($authorXml) | Out-Null
